$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1927
$wsExpo.Range("F4").Value = 830
$wsExpo.Range("F5").Value = 898
$wsExpo.Range("F6").Value = 274

# Sheet "全部类型" - same underlying events, mirrored rows (offset by 1 due to extra row)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1927
$wsAll.Range("F5").Value = 830
$wsAll.Range("F6").Value = 898
$wsAll.Range("F7").Value = 274
